# Add a "team record" (Wins / Losses / Ties) to the player data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row/column of the existing data.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# New columns go immediately after the last existing column (AC -> AD, AE, AF).
$winsCol   = $lastCol + 1
$lossesCol = $lastCol + 2
$tiesCol   = $lastCol + 3

# Header row: new column titles, matching the style of the existing header cells.
$ws.Cells.Item(1, $winsCol).Value   = "Wins"
$ws.Cells.Item(1, $lossesCol).Value = "Losses"
$ws.Cells.Item(1, $tiesCol).Value   = "Ties"

$ws.Cells.Item(1, $lastCol).Copy()
$headerRange = $ws.Range($ws.Cells.Item(1, $winsCol), $ws.Cells.Item(1, $tiesCol))
$headerRange.PasteSpecial(-4122)  # xlPasteFormats

# Every player row gets the same team record: 49 wins, 63 losses, 0 ties.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $winsCol).Value   = 49
    $ws.Cells.Item($r, $lossesCol).Value = 63
    $ws.Cells.Item($r, $tiesCol).Value   = 0
}
